$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Renumber the repeated test rows (previously 4,5,6 -> now 1,2,3) ---
$ws.Range("A9").Value = 1
$ws.Range("A10").Value = 2
$ws.Range("A11").Value = 3

# --- New block: buggfix 2 regression tests (rows 13-15) ---
# Fill column B (descriptions) first, then C (indata), then F (comments),
# matching the order new strings were introduced in the source workbook.
$ws.Range("B13").Value = "Unittest isIsoscelesTest() med En sida = 0"
$ws.Range("B14").Value = "Unittest isScaleneTest()) med En sida = 0"
$ws.Range("B15").Value = "Unittest isEquilateralTest() med En sida = 0"

$ws.Range("C13").Value = "1.2, 0, 1.2"
$ws.Range("C14").Value = "1.2, 0, 1.3"
$ws.Range("C15").Value = "0, 0, 0"

$ws.Range("F13").Value = "Triangel kan ej ha en sida = 0"
$ws.Range("F14").Value = "Triangel kan ej ha en sida = 0"
$ws.Range("F15").Value = "Triangel kan ej ha en sida = 0"

$ws.Range("A13").Value = 4
$ws.Range("D13").Value = "Fail"
$ws.Range("E13").Value = "Pass"

$ws.Range("A14").Value = 5
$ws.Range("D14").Value = "Fail"
$ws.Range("E14").Value = "Pass"

$ws.Range("A15").Value = 6
$ws.Range("D15").Value = "Fail"
$ws.Range("E15").Value = "Pass"

# --- Buggfix 2 code excerpt (rows 17-21) ---
$ws.Range("B17").Value = "Buggfix 2 - kontroll av inparametern i konstruktorn <= 0"
$ws.Range("B18").Value = "if ((a <= 0) || (b <= 0) || (c <= 0))"
$ws.Range("B19").Value = "      {"
$ws.Range("B20").Value = "          throw new ArgumentOutOfRangeException(""Cannot use a value less or equal to 0 as a side in a triangle"");"
$ws.Range("B21").Value = "      }"

# --- Re-run of the three regression tests, before-the-fix results (rows 23-25) ---
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Unittest isIsoscelesTest() med En sida = 0"
$ws.Range("C23").Value = "1.2, 0, 1.2"
$ws.Range("D23").Value = "Fail"
$ws.Range("E23").Value = "Fail"

$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Unittest isScaleneTest()) med En sida = 0"
$ws.Range("C24").Value = "1.2, 0, 1.3"
$ws.Range("D24").Value = "Fail"
$ws.Range("E24").Value = "Fail"

$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "Unittest isEquilateralTest() med En sida = 0"
$ws.Range("C25").Value = "0, 0, 0"
$ws.Range("D25").Value = "Fail"
$ws.Range("E25").Value = "Fail"

# --- Cosmetic sheet/view updates ---
$ws.Columns.Item(2).ColumnWidth = 100.3

$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$ws.Range("F23").Select()
